# Natmi following Dr Hou advice
#
# Rebuilds the LR-pair table (Cx3cl1 -> Itgav) so that the "Sending
# cluster" / "Target cluster" grid covers all three clusters (ECs, FAPs,
# sCs) against each other (a full 3x3 = 9 row block), instead of just
# the FAPs/sCs x ECs/FAPs/sCs 6-row block it used to have, and refreshes
# all of the numeric edge-weight / specificity columns to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (in order): A..D are the cluster / gene labels, E..T are the
# numeric LR-pair statistics. One array per output row, in row order
# starting at row 2.
$data = @(
    @("ECs",  "Cx3cl1", "Itgav", "ECs",  2, 0.6666666666666666, 4.550961,          13.652883,        0.3980483771262702, 0.3980483771262702, 3, 1, 21.09934133333334, 63.29802400000001, 0.2917236204149438, 0.2917236204149438, 96.022279533688,   864.200515803192,   0.1161201136755684,  0.1161201136755684),
    @("ECs",  "Cx3cl1", "Itgav", "FAPs", 2, 0.6666666666666666, 4.550961,          13.652883,        0.3980483771262702, 0.3980483771262702, 3, 1, 35.81943766666667, 107.458313,         0.4952465516465762, 0.4952465516465762, 163.012863862931,  1467.115774766379,  0.1971320861603012,  0.1971320861603012),
    @("ECs",  "Cx3cl1", "Itgav", "sCs",  2, 0.6666666666666666, 4.550961,          13.652883,        0.3980483771262702, 0.3980483771262702, 3, 1, 15.40769666666667, 46.22309,           0.2130298279384801, 0.2130298279384801, 70.11982662983,    631.07843966847,    0.08479617729040057, 0.08479617729040058),
    @("FAPs", "Cx3cl1", "Itgav", "ECs",  3, 1,                  5.449245,          16.347735,        0.4766165055717775, 0.4766165055717775, 3, 1, 21.09934133333334, 63.29802400000001, 0.2917236204149438, 0.2917236204149438, 114.97548026396,   1034.77932237564,   0.1390402925549181,  0.1390402925549182),
    @("FAPs", "Cx3cl1", "Itgav", "FAPs", 3, 1,                  5.449245,          16.347735,        0.4766165055717775, 0.4766165055717775, 3, 1, 35.81943766666667, 107.458313,         0.4952465516465762, 0.4952465516465762, 195.188891607895,  1756.700024471055,  0.236042680842264,   0.236042680842264),
    @("FAPs", "Cx3cl1", "Itgav", "sCs",  3, 1,                  5.449245,          16.347735,        0.4766165055717775, 0.4766165055717775, 3, 1, 15.40769666666667, 46.22309,           0.2130298279384801, 0.2130298279384801, 83.96031402235,    755.64282620115,    0.1015335321745954,  0.1015335321745954),
    @("sCs",  "Cx3cl1", "Itgav", "ECs",  3, 1,                  1.432979666666667, 4.298939,         0.1253351173019523, 0.1253351173019523, 3, 1, 21.09934133333334, 63.29802400000001, 0.2917236204149438, 0.2917236204149438, 30.23492711072622, 272.114343996536,   0.03656321418445719, 0.0365632141844572),
    @("sCs",  "Cx3cl1", "Itgav", "FAPs", 3, 1,                  1.432979666666667, 4.298939,         0.1253351173019523, 0.1253351173019523, 3, 1, 35.81943766666667, 107.458313,         0.4952465516465762, 0.4952465516465762, 51.32852584776744, 461.956732629907,   0.06207178464401101, 0.06207178464401102),
    @("sCs",  "Cx3cl1", "Itgav", "sCs",  3, 1,                  1.432979666666667, 4.298939,         0.1253351173019523, 0.1253351173019523, 3, 1, 15.40769666666667, 46.22309,           0.2130298279384801, 0.2130298279384801, 22.07891603350111, 198.71024430151,    0.02670011847348412, 0.02670011847348412)
)

$r = 2
foreach ($rowVals in $data) {
    $c = 1
    foreach ($v in $rowVals) {
        $ws.Cells.Item($r, $c).Value = $v
        $c++
    }
    $r++
}
